$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has columns: A=code, B=name, C=status,
# D=codeforiati:group-code, E=codeforiati:category-name,
# F=codeforiati:category-code, G=codeforiati:group-name  (BEFORE)
#
# After the edit the header labels for D/E/G are cyclically rotated so that
# D=category-name, E=group-name, F=category-code (unchanged), G=group-code,
# and every data row's values follow the same column rotation so that the
# semantic meaning (category-name/group-name/category-code/group-code) of
# each value stays attached to the correct header.
#
# i.e. for every row: newD = oldE ; newE = oldG ; newF = oldF ; newG = oldD

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rng = $ws.Range("D1:G$lastRow")

# Make sure the values round-trip as TEXT (not auto-converted to numbers),
# matching the original shared-string storage, then restore the default
# "Normal" style so no stray cell formatting is introduced.
$rng.NumberFormat = "@"

$arr = $rng.Value2
$rows = $arr.GetLength(0)
$cols = $arr.GetLength(1)

$new = New-Object 'object[,]' $rows,$cols
for ($i = 1; $i -le $rows; $i++) {
    $oldD = $arr[$i,1]
    $oldE = $arr[$i,2]
    $oldF = $arr[$i,3]
    $oldG = $arr[$i,4]

    $new[$i-1,0] = $oldE
    $new[$i-1,1] = $oldG
    $new[$i-1,2] = $oldF
    $new[$i-1,3] = $oldD
}

$rng.Value2 = $new
$rng.Style = "Normal"
